$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells are stored as literal text (e.g. "283.39", "2.10%", "9"),
# not numbers/percentages. Force Text format on exactly the cells we touch so
# Excel keeps the new values as strings instead of silently parsing them into
# numeric/percentage types when they are written.

# --- Column G (Hora): bump every row from 8 to 9 ---
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "9"

# --- Column D (Price) updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '283.39'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.59'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.065'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06483'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.216'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.425'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9110'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1548'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06545'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07610'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02757'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08964'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001587'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006370'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006075'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.449'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.367'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.966'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.1544'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04456'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001186'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004460'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001199'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0001636'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006623'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1231'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002049'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01243'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005402'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.933'

# --- Column E (Volume(1h)) updates ---
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.10%'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.94%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '3.62%'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.98%'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '3.53%'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '20.44%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.69%'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '27.19%'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.70%'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-4.50%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.12%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.47%'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.33%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.28%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.87%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.59%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.43%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.20%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.82%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.49%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.91%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.66%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.79%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '15.35%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '1.69%'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '-15.70%'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.81%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-2.44%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4.89%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '6.77%'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4.72%'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.70%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.12%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '14.73%'

